$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1. Update the "Date" metadata value on the Metadata sheet (row 8).
# -----------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2026-01-23T08:28:04+00:00"

# -----------------------------------------------------------------
# 2. Insert a new row for "frequenceAdministration" right before the
#    existing "doseAadministrer" row (row 6) on the Elements sheet.
# -----------------------------------------------------------------
$wsEl = $wb.Worksheets.Item("Elements")

# Insert a blank row at position 6 (pushes doseAadministrer and below
# down by one row).
$wsEl.Rows.Item(6).Insert()

# Copy the formatting (border/fill/alignment/style) of the row below
# (now row 7, "doseAadministrer") onto the new blank row 6 so the new
# row matches the rest of the table's look.
$wsEl.Range("A7:AJ7").Copy()
$wsEl.Range("A6:AJ6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's content.
$wsEl.Range("A6").Value = "fr-lm-traitement-subordonne.frequenceAdministration"
$wsEl.Range("B6").Value = "fr-lm-traitement-subordonne.frequenceAdministration"
$wsEl.Range("F6").Value = "0"
$wsEl.Range("G6").Value = "1"
$wsEl.Range("K6").Value = "dateTime
"
$wsEl.Range("L6").Value = "Fréquence d'administration"
$wsEl.Range("M6").Value = "Fréquence d'administration"
$wsEl.Range("AF6").Value = "fr-lm-traitement-subordonne.frequenceAdministration"
$wsEl.Range("AG6").Value = "0"
$wsEl.Range("AH6").Value = "1"

# The remaining "text" columns on the template rows carry an explicit
# empty string value (as opposed to a genuinely blank cell) -- mirror
# that for the new row too.
$emptyStringCols = @("D6","H6","I6","J6","P6","R6","S6","T6","U6","V6","W6","X6","Y6","Z6","AA6","AB6","AC6","AD6","AE6","AI6","AJ6")
foreach ($c in $emptyStringCols) {
    $wsEl.Range($c).Value = ""
}
